$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of cell reference -> new text value, taken verbatim from the
# authoritative diff of the workbook's OOXML.
$changes = [ordered]@{
    D2 = '43.987.38'
    E2 = '  -0.98%  '
    D3 = '2.225.06'
    E3 = '  -0.64%  '
    E4 = '  -1.92%  '
    D5 = '298.37'
    E5 = '  -2.81%  '
    D6 = '90.28'
    E6 = '  -4.14%  '
    D7 = '0.556'
    E7 = '  -2.45%  '
    E8 = '  -0.56%  '
    D9 = '0.492'
    E9 = '  -5.69%  '
    D10 = '32.99'
    E10 = '  -4.76%  '
    E11 = '  -3.05%  '
    E12 = '  -3.27%  '
    E13 = '  -0.62%  '
    D14 = '2.562.60'
    E14 = '  -0.71%  '
    D15 = '2.223.33'
    E15 = '  -0.07%  '
    E16 = '  -1.09%  '
    D17 = '0.776'
    E17 = '  -6.75%  '
    D18 = '43.864.23'
    E18 = '  -0.55%  '
    B19 = 'ShibaInu'
    C19 = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
    D19 = '0.0₃0903'
    E19 = '  -5.21%  '
    B20 = 'InternetComputer(DFINITY)'
    C20 = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
    D20 = '11.67'
    E20 = '  -1.91%  '
    D21 = '5.95'
    E21 = '  -6.03%  '
    D22 = '64.19'
    E22 = '  -1.98%  '
    D23 = '235.51'
    E23 = '  -0.71%  '
    D24 = '2.80'
    E24 = '  -4.84%  '
    E25 = '  -0.37%  '
    E26 = '  -5.86%  '
    B27 = 'InjectiveProtocol'
    C27 = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
    D27 = '39.12'
    E27 = '  +2.57%  '
    B28 = 'Toncoin'
    C28 = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
    D28 = '2.21'
    E28 = '  -0.45%  '
    D29 = '9.35'
    E29 = '  -4.29%  '
    B30 = 'EthereumClassic'
    C30 = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
    D30 = '19.17'
    E30 = '  -3.77%  '
    B31 = 'Monero'
    C31 = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
    D31 = '151.13'
    E31 = '  -1.36%  '
    E32 = '  -8.99%  '
    D33 = '0.0763'
    E33 = '  -3.82%  '
    E34 = '  -6.01%  '
    E35 = '  -1.73%  '
    D36 = '2.83'
    E36 = '  -7.95%  '
    E37 = '  -6.94%  '
    D38 = '1.67'
    E38 = '  -6.70%  '
    D39 = '0.0301'
    E39 = '  +1.43%  '
    B40 = 'NEARProtocol'
    C40 = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
    D40 = '3.16'
    E40 = '  -6.45%  '
    B41 = 'RenderToken'
    C41 = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
    D41 = '3.59'
    E41 = '  -3.94%  '
    E42 = '  -10.90%  '
    E43 = '  -0.83%  '
    D44 = '1.803.89'
    E44 = '  +0.74%  '
    D45 = '1.79'
    E45 = '  +13.30%  '
    D46 = '0.183'
    E46 = '  -3.94%  '
    D47 = '67.78'
    E47 = '  -3.35%  '
    D48 = '94.36'
    E48 = '  -4.12%  '
    E49 = '  -2.94%  '
    D50 = '72.70'
    E50 = '  -7.54%  '
    D51 = '4.60'
    E51 = '  -5.79%  '
}

foreach ($ref in $changes.Keys) {
    $col = ($ref -replace '[0-9]+$', '')
    $newValue = $changes[$ref]
    $range = $ws.Range($ref)
    if ($col -eq 'D') {
        # Column D holds price figures that look like plain numbers
        # (e.g. "298.37") or like multi-dot grouped numbers
        # (e.g. "43.987.38"). Excel's Value setter auto-converts
        # anything that parses as a number, which would silently
        # rewrite the text ("1.00" -> 1) and change the cell type
        # from string to numeric. Force text interpretation while
        # writing, then drop the temporary formatting again so the
        # cell keeps using the workbook's default (General) style.
        $range.NumberFormat = '@'
        $range.Value = $newValue
        $range.ClearFormats()
    } else {
        $range.Value = $newValue
    }
}
